$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format temporarily so numeric-looking strings
# (e.g. "314.28", "27.363.37") are stored as text, matching the source data,
# then clear the temporary formatting so no stray style is left behind.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.363.37'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '1.856.99'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  -0.72%  '
$ws.Range("D5").Value = '314.28'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  -0.99%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.07326'
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("D10").Value = '0.8822'
$ws.Range("E10").Value = '  +0.95%  '
$ws.Range("D11").Value = '19.94'
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").Value = '0.07804'
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("D13").Value = '1.908.98'
$ws.Range("E13").Value = '  +4.12%  '
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").Value = '6.549'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").Value = '91.88'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").Value = '0.000009052'
$ws.Range("E18").Value = '  +1.85%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = '27.376.15'
$ws.Range("E21").Value = '  +2.02%  '
$ws.Range("D22").Value = '5.129'
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("D23").Value = '10.53'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").Value = '2.104.60'
$ws.Range("E24").Value = '  +5.32%  '
$ws.Range("D25").Value = '1.927'
$ws.Range("E25").Value = '  +5.09%  '
$ws.Range("D26").Value = '152.16'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '18.36'
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '2.076'
$ws.Range("E28").Value = '  -1.18%  '
$ws.Range("D29").Value = '5.106'
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").Value = '116.01'
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("D31").Value = '0.08862'
$ws.Range("D32").Value = '0.7711'
$ws.Range("E32").Value = '  +5.71%  '
$ws.Range("D33").Value = '3.041'
$ws.Range("E33").Value = '  +2.20%  '
$ws.Range("E34").Value = '  +3.46%  '
$ws.Range("D35").Value = '4.499'
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").Value = '2.664'
$ws.Range("E36").Value = '  +5.49%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.01959'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.076'
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = '0.05232'
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").Value = '2.955'
$ws.Range("E40").Value = '  +0.97%  '
$ws.Range("D41").Value = '7.021'
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("D42").Value = '0.5140'
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").Value = '8.440'
$ws.Range("E44").Value = '  +2.54%  '
$ws.Range("D45").Value = '0.4820'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("D46").Value = '10.31'
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").Value = '103.03'
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("D49").Value = '1.652'
$ws.Range("E49").Value = '  +1.62%  '
$ws.Range("D50").Value = '0.06222'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").Value = '65.77'
$ws.Range("E51").Value = '  +2.36%  '

$priceRange.ClearFormats()
